$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start from a clean slate for formatting so the leftover "alternating
# medium/thin border + grey fill" look from the old table header is gone.
$ws.Cells.ClearFormats()

# ---------------------------------------------------------------------
# New catalogue rows (11-18). Shared-string pool order matters for an
# exact match, so B16 ("SEPARADA") is written before B15 ("UNIÓN LIBRE")
# to reproduce the original authoring order.
# ---------------------------------------------------------------------
$ws.Range("A11").Value = 11
$ws.Range("B11").Value = "CASADA"

$ws.Range("A12").Value = 12
$ws.Range("B12").Value = "SOLTERA"

$ws.Range("A13").Value = 13
$ws.Range("B13").Value = "DIVORCIADA"

$ws.Range("A14").Value = 14
$ws.Range("B14").Value = "VIUDA"

$ws.Range("A16").Value = 16
$ws.Range("B16").Value = "SEPARADA"

$ws.Range("A15").Value = 15
$ws.Range("B15").Value = "UNIÓN LIBRE"

$ws.Range("A17").Value = 88
$ws.Range("B17").Value = "NO ESPECIFICADO"

$ws.Range("A18").Value = 99
$ws.Range("B18").Value = "SE IGNORA"

# ---------------------------------------------------------------------
# Formatting: a single thin black border around every cell of the table.
# ---------------------------------------------------------------------
$ws.Range("A1:B18").Borders.LineStyle = 1

# Centered, vertically centered, wrapped text for the header and the
# original (0-9) catalogue rows, plus the two new "no especificado /
# se ignora" rows at the bottom.
$ws.Range("A1:B10").HorizontalAlignment = -4108
$ws.Range("A1:B10").VerticalAlignment = -4108
$ws.Range("A1:B10").WrapText = $true

$ws.Range("B17:B18").HorizontalAlignment = -4108
$ws.Range("B17:B18").VerticalAlignment = -4108
$ws.Range("B17:B18").WrapText = $true

# Bold header row.
$ws.Range("A1:B1").Font.Bold = $true

# Leave the selection the way the saved workbook has it.
$ws.Range("B17").Select() | Out-Null
